$d = $word.ActiveDocument

# Paragraph "UPDATE tablo_adı SET alan1=değer1, alan2=değer2 WHERE koşul_ifadesi"
$pUpdate = $d.Paragraphs(36)

# Merge the "UPDATE" run with the following " " run -> "UPDATE "
$pUpdate.Range.Find.Execute("UPDATE ", $true, $false, $false, $false, $false,
                             $true, 1, $false, "UPDATE ", 2)

# Merge the " " run (after "tablo_adı") with "SET alan1=değer1, alan2=değer2 WHERE "
$pUpdate.Range.Find.Execute(" SET alan1=değer1, alan2=değer2 WHERE ", $true, $false, $false, $false, $false,
                             $true, 1, $false, " SET alan1=değer1, alan2=değer2 WHERE ", 2)

# Paragraph "DELETE FROM tablo_adı WHERE koşul_ifadesi"
$pDelete = $d.Paragraphs(42)

# Merge the "DELETE FROM" run with the following " " run -> "DELETE FROM "
$pDelete.Range.Find.Execute("DELETE FROM ", $true, $false, $false, $false, $false,
                             $true, 1, $false, "DELETE FROM ", 2)

# Paragraph containing the TC Kimlik number example - bump 12345678905 -> 12345678906
$pTc = $d.Paragraphs(46)
$pTc.Range.Find.Execute("12345678905", $true, $false, $false, $false, $false,
                         $true, 1, $false, "12345678906", 2)
